$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H129").Value = 831.92
$ws.Range("J129").Value = 904.62067
$ws.Range("L129").Value = 2713.86201
$ws.Range("N129").Value = -12713.86201
$ws.Range("H137").Value = 1454.6111
$ws.Range("I137").Value = 1186.4688
$ws.Range("J137").Value = 3599.75
$ws.Range("K137").Value = 3559.4064
$ws.Range("L137").Value = 10799.25
$ws.Range("M137").Value = -1009.4064
$ws.Range("N137").Value = -15899.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8698.444
$ws.Range("I45").Value = 14084.6
$ws.Range("J45").Value = 1965.75
$ws.Range("K45").Value = 14084.6
$ws.Range("L45").Value = 1965.75
$ws.Range("M45").Value = -13707.6
$ws.Range("N45").Value = -2719.75
$ws.Range("H61").Value = 6891.3335
$ws.Range("I61").Value = 8023.2354
$ws.Range("K61").Value = 8023.2354
$ws.Range("M61").Value = -7811.2354
$ws.Range("H74").Value = 1570.7778
$ws.Range("I74").Value = 1484.84
$ws.Range("J74").Value = 1766.091
$ws.Range("K74").Value = 1484.84
$ws.Range("L74").Value = 1766.091
$ws.Range("M74").Value = -610.8399999999999
$ws.Range("N74").Value = -3514.091
$ws.Range("H77").Value = 1570.7778
$ws.Range("I77").Value = 1484.84
$ws.Range("J77").Value = 1766.091
$ws.Range("K77").Value = 7424.2
$ws.Range("L77").Value = 8830.455
$ws.Range("M77").Value = -3056.2
$ws.Range("N77").Value = -17566.455
$ws.Range("H122").Value = 1426070.6
$ws.Range("I122").Value = 1509780.8
$ws.Range("K122").Value = 4529342.4
$ws.Range("M122").Value = -4526892.4
$ws.Range("H123").Value = 30424
$ws.Range("J123").Value = 30424
$ws.Range("L123").Value = 30424
$ws.Range("N123").Value = -40224
$ws.Range("H132").Value = 5412.0293
$ws.Range("I132").Value = 1453.4375
$ws.Range("J132").Value = 8930.777
$ws.Range("K132").Value = 4360.3125
$ws.Range("L132").Value = 26792.331
$ws.Range("M132").Value = -1830.3125
$ws.Range("N132").Value = -31852.331
$ws.Range("H136").Value = 6891.3335
$ws.Range("I136").Value = 8023.2354
$ws.Range("K136").Value = 24069.7062
$ws.Range("M136").Value = -21519.7062

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4264.75
$ws.Range("I134").Value = 4683.242
$ws.Range("K134").Value = 14049.726
$ws.Range("M134").Value = -11514.726

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4882.1465
$ws.Range("I31").Value = 1200.76
$ws.Range("J31").Value = 10634.3125
$ws.Range("K31").Value = 1200.76
$ws.Range("L31").Value = 10634.3125
$ws.Range("M31").Value = -905.76
$ws.Range("N31").Value = -11224.3125
$ws.Range("H34").Value = 4882.1465
$ws.Range("I34").Value = 1200.76
$ws.Range("J34").Value = 10634.3125
$ws.Range("K34").Value = 1200.76
$ws.Range("L34").Value = 10634.3125
$ws.Range("M34").Value = -998.76
$ws.Range("N34").Value = -11038.3125
$ws.Range("H58").Value = 1180.2084
$ws.Range("I58").Value = 995.55
$ws.Range("K58").Value = 995.55
$ws.Range("M58").Value = -792.55
$ws.Range("H132").Value = 2404.6
$ws.Range("I132").Value = 1876.2858
$ws.Range("J132").Value = 3637.3333
$ws.Range("K132").Value = 5628.857400000001
$ws.Range("L132").Value = 10911.9999
$ws.Range("M132").Value = -3098.857400000001
$ws.Range("N132").Value = -15971.9999
$ws.Range("H134").Value = 1589.5625
$ws.Range("I134").Value = 1533
$ws.Range("J134").Value = 1734.1111
$ws.Range("K134").Value = 4599
$ws.Range("L134").Value = 5202.3333
$ws.Range("M134").Value = -2064
$ws.Range("N134").Value = -10272.3333
$ws.Range("H136").Value = 1180.2084
$ws.Range("I136").Value = 995.55
$ws.Range("K136").Value = 2986.65
$ws.Range("M136").Value = -436.6499999999996

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3814303.5
$ws.Range("I122").Value = 4631161.5
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 13893484.5
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -13891034.5
$ws.Range("N122").Value = -11800
$ws.Range("H126").Value = 6042.7827
$ws.Range("I126").Value = 7143.5
$ws.Range("J126").Value = 2080.2
$ws.Range("K126").Value = 21430.5
$ws.Range("L126").Value = 6240.599999999999
$ws.Range("M126").Value = -18960.5
$ws.Range("N126").Value = -11180.6
$ws.Range("H132").Value = 2935.8684
$ws.Range("I132").Value = 1920.8
$ws.Range("K132").Value = 5762.4
$ws.Range("M132").Value = -3232.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3703316
$ws.Range("I122").Value = 4203693
$ws.Range("J122").Value = 2002034
$ws.Range("K122").Value = 12611079
$ws.Range("L122").Value = 6006102
$ws.Range("M122").Value = -12608629
$ws.Range("N122").Value = -6011002
$ws.Range("H132").Value = 28898242
$ws.Range("I132").Value = 39404884
$ws.Range("J132").Value = 4974.5
$ws.Range("K132").Value = 118214652
$ws.Range("L132").Value = 14923.5
$ws.Range("M132").Value = -118212122
$ws.Range("N132").Value = -19983.5
$ws.Range("H136").Value = 5713.0713
$ws.Range("I136").Value = 5026.3887
$ws.Range("J136").Value = 9833.166999999999
$ws.Range("K136").Value = 15079.1661
$ws.Range("L136").Value = 29499.501
$ws.Range("M136").Value = -12529.1661
$ws.Range("N136").Value = -34599.501

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1887.875
$ws.Range("I122").Value = 1887.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5663.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3213.625
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 1494.5714
$ws.Range("I132").Value = 1030.6923
$ws.Range("K132").Value = 3092.0769
$ws.Range("M132").Value = -562.0769
$ws.Range("H136").Value = 4225.353
$ws.Range("I136").Value = 4960.9165
$ws.Range("J136").Value = 2460
$ws.Range("K136").Value = 14882.7495
$ws.Range("L136").Value = 7380
$ws.Range("M136").Value = -12332.7495
$ws.Range("N136").Value = -12480
